# 6.a.1.1 workbook update:
#  - extend the year series from 2019 through 2023 (4 new columns: X, Y, Z, AA)
#  - hide the old (2000-2007) year columns D:K
#  - bump the row height of the header/data rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 4: years) and the three data rows (5,6,7) ---
# Clone formatting from the last existing year column (W) into the four new
# year columns so the new cells inherit the same borders/alignment/number
# format as the rest of the table, then overwrite with the real values.
$ws.Range("W4").Copy($ws.Range("X4"))
$ws.Range("W4").Copy($ws.Range("Y4"))
$ws.Range("W4").Copy($ws.Range("Z4"))
$ws.Range("W4").Copy($ws.Range("AA4"))

$ws.Range("W5").Copy($ws.Range("X5"))
$ws.Range("W5").Copy($ws.Range("Y5"))
$ws.Range("W5").Copy($ws.Range("Z5"))
$ws.Range("W5").Copy($ws.Range("AA5"))

$ws.Range("W6").Copy($ws.Range("X6"))
$ws.Range("W6").Copy($ws.Range("Y6"))
$ws.Range("W6").Copy($ws.Range("Z6"))
$ws.Range("W6").Copy($ws.Range("AA6"))

$ws.Range("W7").Copy($ws.Range("X7"))
$ws.Range("W7").Copy($ws.Range("Y7"))
$ws.Range("W7").Copy($ws.Range("Z7"))
$ws.Range("W7").Copy($ws.Range("AA7"))

# Row 4 - new year headers
$ws.Range("X4").Value = 2020
$ws.Range("Y4").Value = 2021
$ws.Range("Z4").Value = 2022
$ws.Range("AA4").Value = 2023

# Row 5 - Investment loan
$ws.Range("X5").Value = 23780
$ws.Range("Y5").Value = 44660
$ws.Range("Z5").Value = 25000
$ws.Range("AA5").Value = 13010

# Row 6 - Investment grant
$ws.Range("X6").Value = 38240
$ws.Range("Y6").Value = 7950
$ws.Range("Z6").Value = 23000
$ws.Range("AA6").Value = 16390

# Row 7 - Investment loan and grant
$ws.Range("X7").Value = 62020
$ws.Range("Y7").Value = 52610
$ws.Range("Z7").Value = 48000
$ws.Range("AA7").Value = 29400

# --- Row heights (header band grew a bit taller) ---
$ws.Rows("4").RowHeight = 16.5
$ws.Rows("5").RowHeight = 16.5
$ws.Rows("6").RowHeight = 16.5
$ws.Rows("7").RowHeight = 16.5

# --- Hide the oldest year columns (2000-2007 = D:K) ---
$ws.Range("D1:K1").EntireColumn.ColumnWidth = 0
$ws.Range("D1:K1").EntireColumn.Hidden = $True
